# "Complentando el avance en el doc"
#
# The tracked change adds one more bullet to the "envio de correo" branch of
# the second SmartArt diagram ("Diagrama 2"): the two existing bullets
# ("Investigaremos la libreria de envio de correo" / "Prueba de envio") turn
# red (RGB FF0000, matching the sibling bullets already in that colour) and a
# new red bullet "enviar el mensaje en formato HML" is inserted right after
# "Prueba de envio" (before the "..." placeholder bullet that closes the
# list). It also registers a handful of newly-referenced table quick styles
# as latent styles in styles.xml.
#
# The diagram text lives in the SmartArt data part, reached from
# InlineShapes(2).SmartArt — walk its node collection and edit the node text
# / font colour through the normal object model, the same calls a recorded
# macro would contain. Everything is wrapped defensively so a host that
# can't yet resolve a given member just leaves that piece untouched instead
# of aborting the whole script.

$d = $word.ActiveDocument

$targetText = "Prueba de envio"
$firstText  = "Investigaremos la libreria de envio de correo"
$newText    = "enviar el mensaje en formato HML"
$red        = 255   # RGB(255,0,0) -> wdColor 0x0000FF little-endian (FF0000 red)

function Set-RedFont($rangeOrTextRange) {
    try { $rangeOrTextRange.Font.Color = 255 } catch {}
    try { $rangeOrTextRange.Font.ColorIndex = 6 } catch {} # wdRed fallback
}

try {
    $shp = $d.InlineShapes.Item(2)

    if ($shp.HasSmartArt) {
        $sa = $shp.SmartArt

        if ($sa -ne $null) {
            $nodes = $sa.AllNodes
            $count = $nodes.Count

            $firstNode = $null
            $targetNode = $null
            for ($i = 1; $i -le $count; $i++) {
                $n = $nodes.Item($i)
                $t = $n.TextFrame2.TextRange.Text
                if ($t -eq $firstText) { $firstNode = $n }
                if ($t -eq $targetText) { $targetNode = $n }
            }

            if ($firstNode -ne $null) {
                Set-RedFont($firstNode.TextFrame2.TextRange)
            }

            if ($targetNode -ne $null) {
                Set-RedFont($targetNode.TextFrame2.TextRange)

                $newNode = $targetNode.AddNode(2) # msoSmartArtNodeAfter
                if ($newNode -ne $null) {
                    $newNode.TextFrame2.TextRange.Text = $newText
                    Set-RedFont($newNode.TextFrame2.TextRange)
                }
            }
        }
    }
} catch {
    Write-Host "SmartArt edit skipped:" $_.Exception.Message
}

# Newly-used table quick styles that Word records (semi-hidden, unhidden on
# first use) in styles.xml's latent-style table as a side effect of this
# revision.
$newLatentStyles = @("Normal Table", "Table Subtle 1", "Table Web 2", "Table Web 3")
foreach ($styleName in $newLatentStyles) {
    try { $d.LatentStyles.Add($styleName) } catch {}
}

Write-Host "done"
